# Adding Dragosel to the forecast portfolio
# This shifts the Entsoe "Actual Consumption" data series forward by one day:
# the data that previously represented 45889 (with a trailing row for the last
# interval) now represents 45890, and the table shrinks from 40 to 39 data rows
# (A1:B41 -> A1:B40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; A=5371; B=45890},
    @{Row=3; A=5337; B=45890.01041666666},
    @{Row=4; A=5300; B=45890.02083333334},
    @{Row=5; A=5252; B=45890.03125},
    @{Row=6; A=5213; B=45890.04166666666},
    @{Row=7; A=5190; B=45890.05208333334},
    @{Row=8; A=5171; B=45890.0625},
    @{Row=9; A=5136; B=45890.07291666666},
    @{Row=10; A=5121; B=45890.08333333334},
    @{Row=11; A=5098; B=45890.09375},
    @{Row=12; A=5085; B=45890.10416666666},
    @{Row=13; A=5085; B=45890.11458333334},
    @{Row=14; A=5083; B=45890.125},
    @{Row=15; A=5099; B=45890.13541666666},
    @{Row=16; A=5114; B=45890.14583333334},
    @{Row=17; A=5123; B=45890.15625},
    @{Row=18; A=5174; B=45890.16666666666},
    @{Row=19; A=5213; B=45890.17708333334},
    @{Row=20; A=5261; B=45890.1875},
    @{Row=21; A=5312; B=45890.19791666666},
    @{Row=22; A=5455; B=45890.20833333334},
    @{Row=23; A=5487; B=45890.21875},
    @{Row=24; A=5516; B=45890.22916666666},
    @{Row=25; A=5585; B=45890.23958333334},
    @{Row=26; A=5792; B=45890.25},
    @{Row=27; A=5837; B=45890.26041666666},
    @{Row=28; A=5886; B=45890.27083333334},
    @{Row=29; A=5953; B=45890.28125},
    @{Row=30; A=6013; B=45890.29166666666},
    @{Row=31; A=5994; B=45890.30208333334},
    @{Row=32; A=5974; B=45890.3125},
    @{Row=33; A=5939; B=45890.32291666666},
    @{Row=34; A=5814; B=45890.33333333334},
    @{Row=35; A=5851; B=45890.34375},
    @{Row=36; A=5812; B=45890.35416666666},
    @{Row=37; A=5732; B=45890.36458333334},
    @{Row=38; A=5712; B=45890.375},
    @{Row=39; A=5641; B=45890.38541666666},
    @{Row=40; A=5623; B=45890.39583333334}
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 1).Value = $entry.A
    $ws.Cells.Item($entry.Row, 2).Value = $entry.B
}

# Remove the now-obsolete last row (previously row 41), shrinking the table
# from 40 data rows to 39 data rows.
$ws.Rows.Item(41).Delete()
